$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Cells.Item(2, 4).Value = "'29.271.54"
$ws.Cells.Item(2, 5).Value = "  -0.33%  "
$ws.Cells.Item(3, 4).Value = "'1.864.78"
$ws.Cells.Item(3, 5).Value = "  -0.48%  "
$ws.Cells.Item(4, 5).Value = "  +0.03%  "
$ws.Cells.Item(5, 4).Value = "'0.7164"
$ws.Cells.Item(5, 5).Value = "  +0.58%  "
$ws.Cells.Item(6, 4).Value = "'240.88"
$ws.Cells.Item(7, 5).Value = "  +0.05%  "
$ws.Cells.Item(8, 4).Value = "'0.3098"
$ws.Cells.Item(8, 5).Value = "  +0.69%  "
$ws.Cells.Item(9, 4).Value = "'0.07719"
$ws.Cells.Item(9, 5).Value = "  -0.90%  "
$ws.Cells.Item(10, 4).Value = "'25.09"
$ws.Cells.Item(10, 5).Value = "  +1.07%  "
$ws.Cells.Item(11, 2).Value = "TRON"
$ws.Cells.Item(11, 3).Value = "https://coinranking.com/coin/qUhEFk1I61atv+tron-trx"
$ws.Cells.Item(11, 4).Value = "'0.08317"
$ws.Cells.Item(11, 5).Value = "  +0.96%  "
$ws.Cells.Item(12, 2).Value = "WrappedEther"
$ws.Cells.Item(12, 3).Value = "https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth"
$ws.Cells.Item(12, 4).Value = "'1.975.00"
$ws.Cells.Item(12, 5).Value = "  +5.62%  "
$ws.Cells.Item(13, 4).Value = "'0.7184"
$ws.Cells.Item(13, 5).Value = "  -0.63%  "
$ws.Cells.Item(14, 4).Value = "'5.224"
$ws.Cells.Item(14, 5).Value = "  -0.72%  "
$ws.Cells.Item(15, 4).Value = "'90.99"
$ws.Cells.Item(15, 5).Value = "  -0.50%  "
$ws.Cells.Item(16, 4).Value = "'29.416.89"
$ws.Cells.Item(16, 5).Value = "  -0.29%  "
$ws.Cells.Item(17, 4).Value = "'5.954"
$ws.Cells.Item(17, 5).Value = "  +1.20%  "
$ws.Cells.Item(18, 4).Value = "'243.94"
$ws.Cells.Item(18, 5).Value = "  +0.16%  "
$ws.Cells.Item(19, 2).Value = "WrappedliquidstakedEther2.0"
$ws.Cells.Item(19, 3).Value = "https://coinranking.com/coin/CiixT63n3+wrappedliquidstakedether20-wsteth"
$ws.Cells.Item(19, 4).Value = "'2.168.32"
$ws.Cells.Item(19, 5).Value = "  +2.28%  "
$ws.Cells.Item(20, 2).Value = "ShibaInu"
$ws.Cells.Item(20, 3).Value = "https://coinranking.com/coin/xz24e0BjL+shibainu-shib"
$ws.Cells.Item(20, 4).Value = "'0.000007816"
$ws.Cells.Item(20, 5).Value = "  -1.28%  "
$ws.Cells.Item(21, 4).Value = "'13.16"
$ws.Cells.Item(21, 5).Value = "  -0.76%  "
$ws.Cells.Item(22, 2).Value = "Chainlink"
$ws.Cells.Item(22, 3).Value = "https://coinranking.com/coin/VLqpJwogdhHNb+chainlink-link"
$ws.Cells.Item(22, 4).Value = "'8.043"
$ws.Cells.Item(22, 5).Value = "  +0.86%  "
$ws.Cells.Item(23, 2).Value = "Dai"
$ws.Cells.Item(23, 3).Value = "https://coinranking.com/coin/MoTuySvg7+dai-dai"
$ws.Cells.Item(23, 4).Value = "'1.001"
$ws.Cells.Item(23, 5).Value = "  +0.11%  "
$ws.Cells.Item(24, 5).Value = "  -0.01%  "
$ws.Cells.Item(25, 4).Value = "'0.1618"
$ws.Cells.Item(25, 5).Value = "  +4.87%  "
$ws.Cells.Item(26, 4).Value = "'162.91"
$ws.Cells.Item(26, 5).Value = "  -0.31%  "
$ws.Cells.Item(27, 4).Value = "'8.935"
$ws.Cells.Item(27, 5).Value = "  -0.49%  "
$ws.Cells.Item(28, 4).Value = "'18.63"
$ws.Cells.Item(28, 5).Value = "  +1.82%  "
$ws.Cells.Item(29, 4).Value = "'1.345"
$ws.Cells.Item(29, 5).Value = "  -0.90%  "
$ws.Cells.Item(30, 2).Value = "Filecoin"
$ws.Cells.Item(30, 3).Value = "https://coinranking.com/coin/ymQub4fuB+filecoin-fil"
$ws.Cells.Item(30, 4).Value = "'4.455"
$ws.Cells.Item(30, 5).Value = "  +1.99%  "
$ws.Cells.Item(31, 2).Value = "PancakeSwap"
$ws.Cells.Item(31, 3).Value = "https://coinranking.com/coin/ncYFcP709+pancakeswap-cake"
$ws.Cells.Item(31, 4).Value = "'1.497"
$ws.Cells.Item(31, 5).Value = "  +0.77%  "
$ws.Cells.Item(32, 4).Value = "'4.249"
$ws.Cells.Item(32, 5).Value = "  +3.44%  "
$ws.Cells.Item(33, 4).Value = "'0.05191"
$ws.Cells.Item(33, 5).Value = "  -1.44%  "
$ws.Cells.Item(34, 5).Value = "  +0.50%  "
$ws.Cells.Item(35, 5).Value = "  -1.88%  "
$ws.Cells.Item(36, 4).Value = "'0.7276"
$ws.Cells.Item(36, 5).Value = "  +1.34%  "
$ws.Cells.Item(37, 4).Value = "'2.681"
$ws.Cells.Item(37, 5).Value = "  +0.07%  "
$ws.Cells.Item(38, 4).Value = "'0.01861"
$ws.Cells.Item(38, 5).Value = "  +0.14%  "
$ws.Cells.Item(39, 4).Value = "'2.699"
$ws.Cells.Item(39, 5).Value = "  -0.34%  "
$ws.Cells.Item(40, 4).Value = "'1.182.72"
$ws.Cells.Item(40, 5).Value = "  -2.05%  "
$ws.Cells.Item(41, 2).Value = "FraxShare"
$ws.Cells.Item(41, 3).Value = "https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs"
$ws.Cells.Item(41, 4).Value = "'6.238"
$ws.Cells.Item(41, 5).Value = "  +3.03%  "
$ws.Cells.Item(42, 2).Value = "TrustWalletToken"
$ws.Cells.Item(42, 3).Value = "https://coinranking.com/coin/Hm3OlynlC+trustwallettoken-twt"
$ws.Cells.Item(42, 4).Value = "'0.9048"
$ws.Cells.Item(42, 5).Value = "  -0.14%  "
$ws.Cells.Item(43, 4).Value = "'73.19"
$ws.Cells.Item(43, 5).Value = "  +1.28%  "
$ws.Cells.Item(44, 4).Value = "'1.001"
$ws.Cells.Item(44, 5).Value = "  +0.02%  "
$ws.Cells.Item(45, 4).Value = "'102.89"
$ws.Cells.Item(45, 5).Value = "  -0.35%  "
$ws.Cells.Item(46, 4).Value = "'2.059.86"
$ws.Cells.Item(46, 5).Value = "  +2.20%  "
$ws.Cells.Item(47, 4).Value = "'0.5221"
$ws.Cells.Item(47, 5).Value = "  -2.25%  "
$ws.Cells.Item(48, 4).Value = "'1.778"
$ws.Cells.Item(48, 5).Value = "  +1.38%  "
$ws.Cells.Item(49, 2).Value = "EnergySwap"
$ws.Cells.Item(49, 3).Value = "https://coinranking.com/coin/SbWqqTui-+energyswap-ens"
$ws.Cells.Item(49, 4).Value = "'9.359"
$ws.Cells.Item(49, 5).Value = "  +1.45%  "
$ws.Cells.Item(50, 2).Value = "Frax"
$ws.Cells.Item(50, 3).Value = "https://coinranking.com/coin/KfWtaeV1W+frax-frax"
$ws.Cells.Item(50, 4).Value = "'1.021"
$ws.Cells.Item(50, 5).Value = "  +2.11%  "
$ws.Cells.Item(51, 2).Value = "SynthetixNetwork"
$ws.Cells.Item(51, 3).Value = "https://coinranking.com/coin/sgxZRXbK0FDc+synthetixnetwork-snx"
$ws.Cells.Item(51, 4).Value = "'2.866"
$ws.Cells.Item(51, 5).Value = "  -0.97%  "
